$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B136").Value = 63902
$ws.Range("D136").Value = 32.02
$ws.Range("E136").Value = 34.04
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0

$ws.Range("B137").Value = 48654
$ws.Range("D137").Value = 32.02
$ws.Range("E137").Value = 38.26
$ws.Range("F137").Value = -1
$ws.Range("G137").Value = -32.02

$ws.Range("B146").Value = 64350
$ws.Range("D146").Value = 66.44
$ws.Range("E146").Value = 70.63
$ws.Range("F146").Value = 2
$ws.Range("G146").Value = 132.88

$ws.Range("B147").Value = 57756
$ws.Range("D147").Value = 66.44
$ws.Range("E147").Value = 79.37
$ws.Range("F147").Value = -100
$ws.Range("G147").Value = -6644

$ws.Range("B148").Value = 53925
$ws.Range("D148").Value = 66.44
$ws.Range("E148").Value = 79.37
$ws.Range("F148").Value = 1
$ws.Range("G148").Value = 66.44

$ws.Range("B277").Value = 61610
$ws.Range("D277").Value = 102.71
$ws.Range("E277").Value = 122.71
$ws.Range("F277").Value = -58
$ws.Range("G277").Value = -5957.18

$ws.Range("B278").Value = 63565
$ws.Range("D278").Value = 102.71
$ws.Range("E278").Value = 109.19
$ws.Range("F278").Value = 60
$ws.Range("G278").Value = 6162.6

$ws.Range("B295").Value = 63531
$ws.Range("D295").Value = 143.48
$ws.Range("E295").Value = 152.53
$ws.Range("F295").Value = 80
$ws.Range("G295").Value = 11478.4

$ws.Range("B296").Value = 63571
$ws.Range("D296").Value = 143.48
$ws.Range("E296").Value = 152.53
$ws.Range("F296").Value = 4
$ws.Range("G296").Value = 573.92

$ws.Range("B356").Value = 63681
$ws.Range("D356").Value = 22.42
$ws.Range("E356").Value = 23.84
$ws.Range("F356").Value = 0
$ws.Range("G356").Value = 0

$ws.Range("B357").Value = 31930
$ws.Range("D357").Value = 22.42
$ws.Range("E357").Value = 26.8
$ws.Range("F357").Value = -62
$ws.Range("G357").Value = -1390.04

$ws.Range("B465").Value = 53757
$ws.Range("D465").Value = 13.45
$ws.Range("E465").Value = 16.08
$ws.Range("F465").Value = -159
$ws.Range("G465").Value = -2138.55

$ws.Range("B466").Value = 65069
$ws.Range("D466").Value = 13.45
$ws.Range("E466").Value = 14.3
$ws.Range("F466").Value = 2
$ws.Range("G466").Value = 26.9

$ws.Range("B472").Value = 45695
$ws.Range("D472").Value = 19.73
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28

$ws.Range("B473").Value = 64915
$ws.Range("D473").Value = 19.73
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0

$ws.Range("B479").Value = 64927
$ws.Range("D479").Value = 16.22
$ws.Range("E479").Value = 17.26
$ws.Range("F479").Value = 168
$ws.Range("G479").Value = 2724.96

$ws.Range("B480").Value = 45718
$ws.Range("D480").Value = 16.22
$ws.Range("E480").Value = 19.38
$ws.Range("F480").Value = -294
$ws.Range("G480").Value = -4768.68

$ws.Range("B485").Value = 64925
$ws.Range("D485").Value = 13.15
$ws.Range("E485").Value = 13.97
$ws.Range("F485").Value = 190
$ws.Range("G485").Value = 2498.5

$ws.Range("B486").Value = 45709
$ws.Range("D486").Value = 13.15
$ws.Range("E486").Value = 15.69
$ws.Range("F486").Value = -300
$ws.Range("G486").Value = -3945

$ws.Range("B487").Value = 45702
$ws.Range("D487").Value = 26.3
$ws.Range("E487").Value = 31.43
$ws.Range("F487").Value = -215
$ws.Range("G487").Value = -5654.5

$ws.Range("B488").Value = 64919
$ws.Range("D488").Value = 26.3
$ws.Range("E488").Value = 27.97
$ws.Range("F488").Value = 94
$ws.Range("G488").Value = 2472.2

$ws.Range("B574").Value = 64810
$ws.Range("D574").Value = 273.92
$ws.Range("E574").Value = 291.22
$ws.Range("F574").Value = 6
$ws.Range("G574").Value = 1643.52

$ws.Range("B575").Value = 53319
$ws.Range("D575").Value = 273.92
$ws.Range("E575").Value = 310.64
$ws.Range("F575").Value = -6
$ws.Range("G575").Value = -1643.52

$ws.Range("B606").Value = 64830
$ws.Range("D606").Value = 32.83
$ws.Range("E606").Value = 34.9
$ws.Range("F606").Value = 112
$ws.Range("G606").Value = 3676.96

$ws.Range("B607").Value = 60022
$ws.Range("D607").Value = 32.83
$ws.Range("E607").Value = 37.22
$ws.Range("F607").Value = -113
$ws.Range("G607").Value = -3709.79

$ws.Range("B715").Value = 63150
$ws.Range("D715").Value = 75.68000000000001
$ws.Range("E715").Value = 80.45
$ws.Range("F715").Value = 55
$ws.Range("G715").Value = 4162.4

$ws.Range("B716").Value = 61428
$ws.Range("D716").Value = 69.16
$ws.Range("E716").Value = 73.52
$ws.Range("F716").Value = 1
$ws.Range("G716").Value = 69.16

$ws.Range("B742").Value = 65079
$ws.Range("D742").Value = 40.87
$ws.Range("E742").Value = 43.44
$ws.Range("F742").Value = 21
$ws.Range("G742").Value = 858.27

$ws.Range("B743").Value = 65362
$ws.Range("D743").Value = 40.87
$ws.Range("E743").Value = 43.44
$ws.Range("F743").Value = 49
$ws.Range("G743").Value = 2002.63

Write-Host "Edit complete"
